# Student Management - Academic Year Code Implementation
# Updates the single approved-history row: corrected student name, a new
# approval date, a fuller "Reason" note, and the person who actually
# picked the student up - then widens the Reason/PickedUp columns so the
# longer text is readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the stray "1" typo in the student's name.
$ws.Range("A2").Value = "John A. Doe"

# Update the approval date.
$ws.Range("E2").Value = "26-03-2025"

# Replace the vague reason with a clear explanation, and record who
# actually picked the student up.
$ws.Range("H2").Value = "Student requires early pickup due to health concerns"
$ws.Range("I2").Value = "Meera Kapoor"

# Widen the Reason (H) and PickedUp (I) columns to fit the new text.
$ws.Columns.Item(8).ColumnWidth = 46.3
$ws.Columns.Item(9).ColumnWidth = 12.45
